# Apply updated crypto price/volume figures (and the Celestia/LidoDAOToken row-order
# swap) to the sheet, matching the scraped GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.235.72'
$ws.Range('E2').Value = '  -1.85%  '

$ws.Range('D3').Value = '2.171.34'
$ws.Range('E3').Value = '  -2.01%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = '''236.32'
$ws.Range('E5').Value = '  -2.30%  '

$ws.Range('D6').Value = '''0.616'
$ws.Range('E6').Value = '  -1.56%  '

$ws.Range('D7').Value = '''70.12'
$ws.Range('E7').Value = '  -4.68%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '''0.580'
$ws.Range('E9').Value = '  -5.20%  '

$ws.Range('D10').Value = '''40.26'
$ws.Range('E10').Value = '  -7.81%  '

$ws.Range('D11').Value = '''0.0926'
$ws.Range('E11').Value = '  -3.30%  '

$ws.Range('D12').Value = '''55.06'
$ws.Range('E12').Value = '  -4.21%  '

$ws.Range('D13').Value = '''6.79'
$ws.Range('E13').Value = '  -4.32%  '

$ws.Range('E14').Value = '  -2.24%  '

$ws.Range('D15').Value = '2.497.39'
$ws.Range('E15').Value = '  -2.04%  '

$ws.Range('D16').Value = '''13.91'
$ws.Range('E16').Value = '  -2.36%  '

$ws.Range('D17').Value = '''0.805'
$ws.Range('E17').Value = '  -4.63%  '

$ws.Range('D18').Value = '2.167.51'
$ws.Range('E18').Value = '  -2.96%  '

$ws.Range('D19').Value = '41.068.05'
$ws.Range('E19').Value = '  -1.91%  '

$ws.Range('D20').Value = '''0.0000101'
$ws.Range('E20').Value = '  -6.67%  '

$ws.Range('D21').Value = '''70.48'
$ws.Range('E21').Value = '  -2.89%  '

$ws.Range('D22').Value = '''5.95'
$ws.Range('E22').Value = '  -3.29%  '

$ws.Range('D23').Value = '''9.77'
$ws.Range('E23').Value = '  -7.39%  '

$ws.Range('D24').Value = '''226.70'
$ws.Range('E24').Value = '  -1.23%  '

$ws.Range('D25').Value = '''1.95'
$ws.Range('E25').Value = '  -6.52%  '

$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('D27').Value = '''10.88'
$ws.Range('E27').Value = '  -5.62%  '

$ws.Range('D28').Value = '''3.54'
$ws.Range('E28').Value = '  -1.22%  '

$ws.Range('E29').Value = '  -2.38%  '

$ws.Range('D31').Value = '''167.95'
$ws.Range('E31').Value = '  +0.71%  '

$ws.Range('D32').Value = '''19.98'
$ws.Range('E32').Value = '  -2.93%  '

$ws.Range('E33').Value = '  +6.18%  '

$ws.Range('D34').Value = '''0.0770'
$ws.Range('E34').Value = '  -3.12%  '

$ws.Range('D35').Value = '''5.16'
$ws.Range('E35').Value = '  -7.93%  '

$ws.Range('E36').Value = '  -3.38%  '

$ws.Range('E37').Value = '  -8.92%  '

$ws.Range('D38').Value = '''4.12'
$ws.Range('E38').Value = '  -3.46%  '

$ws.Range('D39').Value = '''0.0285'
$ws.Range('E39').Value = '  -5.52%  '

$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''2.08'
$ws.Range('E40').Value = '  -1.88%  '

$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '''11.91'
$ws.Range('E41').Value = '  -7.30%  '

$ws.Range('D42').Value = '''5.45'
$ws.Range('E42').Value = '  -3.21%  '

$ws.Range('D43').Value = '''59.93'
$ws.Range('E43').Value = '  -9.85%  '

$ws.Range('D44').Value = '''0.191'
$ws.Range('E44').Value = '  -4.78%  '

$ws.Range('D45').Value = '''8.32'
$ws.Range('E45').Value = '  -4.75%  '

$ws.Range('D46').Value = '''0.0977'
$ws.Range('E46').Value = '  -2.93%  '

$ws.Range('D47').Value = '''97.94'
$ws.Range('E47').Value = '  -5.69%  '

$ws.Range('E48').Value = '  -2.17%  '

$ws.Range('D49').Value = '''1.13'
$ws.Range('E49').Value = '  -2.60%  '

$ws.Range('D50').Value = '''2.22'
$ws.Range('E50').Value = '  -8.71%  '

$ws.Range('D51').Value = '''2.63'
$ws.Range('E51').Value = '  -2.87%  '
